# Weekly fruit/vegetable price update:
# Insert a new row at row 92 (shifting the existing rows 92-97 down to 93-98)
# and populate it with the newest week's data, matching the pattern used by
# all the other rows in this price-history sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 92:97 down to 93:98 to make room for the new weekly entry.
$ws.Range("A92:R92").EntireRow.Insert()

# Populate the new row 92 with this week's record.
$ws.Range("A92").Value = 5
$ws.Range("B92").Value = "Macroferia Regional de Talca"
$ws.Range("C92").Value = "Maule"
$ws.Range("D92").Value = 44585
$ws.Range("D92").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E92").Value = 7
$ws.Range("F92").Value = 100112030
$ws.Range("G92").Value = "Poroto granado"
$ws.Range("H92").Value = "Sin especificar"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 400
$ws.Range("K92").Value = 23000
$ws.Range("L92").Value = 23000
$ws.Range("M92").Value = 23000
$ws.Range("N92").Value = "$/saco 25 kilos"
$ws.Range("O92").Value = "Región del Maule"
$ws.Range("P92").Value = 920
$ws.Range("Q92").Value = 25
$ws.Range("R92").Value = "Hortaliza"
